$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before the current "cierre_tpv_desc" column (E), shifting
# E->F, F->G, G->H, H->I. The new column will hold "cierre_tpv_id".
$ws.Columns.Item(5).Insert()

# The inserted column inherits column D's date style; reset it back to the
# default (unstyled / General) cell style before writing the new data.
$ws.Range("E1:E11").Style = "Normal"

# New header for the inserted column.
$ws.Range("E1").Value = "cierre_tpv_id"

# Rewrite every data row (2-11) with the final values, now laid out across
# columns A-I and re-ordered / re-grouped by cierre_tpv_id.
$ws.Range("B2").Value = "V1"
$ws.Range("C2").Value = "SERVIDOR TIENDA"
$ws.Range("D2").Value = 45700
$ws.Range("E2").Value = 8868
$ws.Range("F2").Value = "Mañana"
$ws.Range("G2").Value = "SMS"
$ws.Range("H2").Value = 54.5
$ws.Range("I2").Value = 3

$ws.Range("B3").Value = "V2"
$ws.Range("C3").Value = "BAR"
$ws.Range("D3").Value = 45700
$ws.Range("E3").Value = 8869
$ws.Range("F3").Value = "Mañana"
$ws.Range("G3").Value = "EUROS"
$ws.Range("H3").Value = 63.5
$ws.Range("I3").Value = 15

$ws.Range("B4").Value = "V2"
$ws.Range("C4").Value = "BAR"
$ws.Range("D4").Value = 45700
$ws.Range("E4").Value = 8869
$ws.Range("F4").Value = "Mañana"
$ws.Range("G4").Value = "TARJETA VISA"
$ws.Range("H4").Value = 173.7
$ws.Range("I4").Value = 32

$ws.Range("B5").Value = "V1"
$ws.Range("C5").Value = "SERVIDOR TIENDA"
$ws.Range("D5").Value = 45700
$ws.Range("E5").Value = 8870
$ws.Range("F5").Value = "Mañana"
$ws.Range("G5").Value = "EUROS"
$ws.Range("H5").Value = 573.8
$ws.Range("I5").Value = 82

$ws.Range("B6").Value = "V1"
$ws.Range("C6").Value = "SERVIDOR TIENDA"
$ws.Range("D6").Value = 45700
$ws.Range("E6").Value = 8870
$ws.Range("F6").Value = "Mañana"
$ws.Range("G6").Value = "TARJETA VISA"
$ws.Range("H6").Value = 1134.89
$ws.Range("I6").Value = 151

$ws.Range("B7").Value = "V2"
$ws.Range("C7").Value = "BAR"
$ws.Range("D7").Value = 45700
$ws.Range("E7").Value = 8871
$ws.Range("F7").Value = "Mañana"
$ws.Range("G7").Value = "EUROS"
$ws.Range("H7").Value = 249.3
$ws.Range("I7").Value = 30

$ws.Range("B8").Value = "V2"
$ws.Range("C8").Value = "BAR"
$ws.Range("D8").Value = 45700
$ws.Range("E8").Value = 8871
$ws.Range("F8").Value = "Mañana"
$ws.Range("G8").Value = "TARJETA VISA"
$ws.Range("H8").Value = 370.77
$ws.Range("I8").Value = 48

$ws.Range("B9").Value = "V1"
$ws.Range("C9").Value = "SERVIDOR TIENDA"
$ws.Range("D9").Value = 45700
$ws.Range("E9").Value = 8872
$ws.Range("F9").Value = "Mañana"
$ws.Range("G9").Value = "EUROS"
$ws.Range("H9").Value = 223.78
$ws.Range("I9").Value = 40

$ws.Range("B10").Value = "V1"
$ws.Range("C10").Value = "SERVIDOR TIENDA"
$ws.Range("D10").Value = 45700
$ws.Range("E10").Value = 8872
$ws.Range("F10").Value = "Mañana"
$ws.Range("G10").Value = "TARJETA VISA"
$ws.Range("H10").Value = 837.98
$ws.Range("I10").Value = 87

$ws.Range("B11").Value = "V2"
$ws.Range("C11").Value = "BAR"
$ws.Range("D11").Value = 45700
$ws.Range("E11").Value = 8873
$ws.Range("F11").Value = "Mañana"
$ws.Range("G11").Value = "TARJETA VISA"
$ws.Range("H11").Value = 126.2
$ws.Range("I11").Value = 19
